$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 123 (shifts old rows 123:216 down to 126:219)
$ws.Range("A123:R125").EntireRow.Insert()

# Common / constant values shared with the surrounding rows in this data block
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112008
$categoria = "Coliflor"
$variedad  = "Sin especificar"
$calidad   = "Primera"
$unidad    = "`$/unidad"
$kgOUnid   = 1
$clasif    = "Hortaliza"

# New row 123
$r = 123
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44452
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $catId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = $calidad
$ws.Cells.Item($r, 10).Value = 1500
$ws.Cells.Item($r, 11).Value = 800
$ws.Cells.Item($r, 12).Value = 800
$ws.Cells.Item($r, 13).Value = 800
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($r, 16).Value = 800
$ws.Cells.Item($r, 17).Value = $kgOUnid
$ws.Cells.Item($r, 18).Value = $clasif

# New row 124
$r = 124
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44452
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $catId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = $calidad
$ws.Cells.Item($r, 10).Value = 2000
$ws.Cells.Item($r, 11).Value = 800
$ws.Cells.Item($r, 12).Value = 800
$ws.Cells.Item($r, 13).Value = 800
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 16).Value = 800
$ws.Cells.Item($r, 17).Value = $kgOUnid
$ws.Cells.Item($r, 18).Value = $clasif

# New row 125
$r = 125
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44452
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $catId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = $calidad
$ws.Cells.Item($r, 10).Value = 2000
$ws.Cells.Item($r, 11).Value = 1000
$ws.Cells.Item($r, 12).Value = 1000
$ws.Cells.Item($r, 13).Value = 1000
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región del Maule"
$ws.Cells.Item($r, 16).Value = 1000
$ws.Cells.Item($r, 17).Value = $kgOUnid
$ws.Cells.Item($r, 18).Value = $clasif
